# Natmi following Dr Hou advice
# Reorders sCs/Btc label and recomputes the Btc-Erbb2 LR table, adding the
# FAPs/sCs x FAPs/sCs (excluding FAPs->FAPs... ) combinations -> rows 2-7.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("FAPs","Btc","Erbb2","ECs",3,1,0.637502,1.912506,0.6183852387299037,0.6183852387299037,3,1,3.709791333333333,11.129374,0.4283284425582907,0.4283284425582907,2.364999394582667,21.284994551244,0.2648719862062164,0.2648719862062164),
    @("FAPs","Btc","Erbb2","FAPs",3,1,0.637502,1.912506,0.6183852387299037,0.6183852387299037,3,1,3.599001333333334,10.797004,0.4155367505499981,0.4155367505499982,2.294370548002667,20.649334932024,0.256961792689909,0.256961792689909),
    @("FAPs","Btc","Erbb2","sCs",3,1,0.637502,1.912506,0.6183852387299037,0.6183852387299037,3,1,1.352297666666667,4.056893,0.1561348068917112,0.1561348068917112,0.8620924670953333,7.758832203858,0.09655145983377823,0.09655145983377823),
    @("sCs","Btc","Erbb2","ECs",2,0.6666666666666666,0.393412,1.180236,0.3816147612700962,0.3816147612700962,3,1,3.709791333333333,11.129374,0.4283284425582907,0.4283284425582907,1.459476428029334,13.135287852264,0.1634564563520742,0.1634564563520743),
    @("sCs","Btc","Erbb2","FAPs",2,0.6666666666666666,0.393412,1.180236,0.3816147612700962,0.3816147612700962,3,1,3.599001333333334,10.797004,0.4155367505499981,0.4155367505499982,1.415890312549334,12.743012812944,0.1585749578600891,0.1585749578600891),
    @("sCs","Btc","Erbb2","sCs",2,0.6666666666666666,0.393412,1.180236,0.3816147612700962,0.3816147612700962,3,1,1.352297666666667,4.056893,0.1561348068917112,0.1561348068917112,0.5320101296386667,4.788091166748,0.05958334705793294,0.05958334705793294)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $i + 2
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $row[$j]
    }
}
